# Update the diary workbook:
#  - Extend the entry for "5 joulu" (row 42): change its time cell to a
#    range, expand the description text, and add hours logged (1.5).
#  - Add a new entry (row 43) for "8 joulu".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42: expand existing "5 joulu" entry ---
$ws.Range("B42").Value = "11.45-13.15"
$ws.Range("C42").Value = "Laskentavarjostin scene setup, fluiditutoriaalin läpikatselu ja seuraavien askelmerkkien suunnittelu"
$ws.Range("G42").Value = 1.5

# --- Row 43: new "8 joulu" entry ---
$ws.Range("A43").Value = "8 joulu"
$ws.Range("B43").Value = "17.45-18.30,"
$ws.Range("C43").Value = "Fluiditutoriaalin palastelua"

# Match formatting of the row above for the new cells.
$ws.Range("B43").NumberFormat = $ws.Range("B42").NumberFormat
$ws.Range("C43").WrapText = $ws.Range("C42").WrapText

$ws.Range("B43").Select()
